$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.076.74"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3
$ws.Range("D3").Value = "1.678.23"
$ws.Range("E3").Value = "  +0.38%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'215.15"
$ws.Range("E5").Value = "  +0.12%  "

# Row 6
$ws.Range("E6").Value = "  -0.25%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("E8").Value = "  +1.93%  "

# Row 9
$ws.Range("D9").Value = "'21.27"
$ws.Range("E9").Value = "  +5.48%  "

# Row 10
$ws.Range("E10").Value = "  +0.31%  "

# Row 11
$ws.Range("D11").Value = "'0.0885"
$ws.Range("E11").Value = "  -0.48%  "

# Row 12
$ws.Range("D12").Value = "1.913.88"

# Row 13
$ws.Range("D13").Value = "1.680.69"
$ws.Range("E13").Value = "  +0.10%  "

# Row 14
$ws.Range("E14").Value = "  +1.06%  "

# Row 15
$ws.Range("E15").Value = "  +1.63%  "

# Row 16
$ws.Range("D16").Value = "'66.06"
$ws.Range("E16").Value = "  +0.77%  "

# Row 17
$ws.Range("D17").Value = "27.070.03"
$ws.Range("E17").Value = "  +0.47%  "

# Row 18
$ws.Range("D18").Value = "'237.10"
$ws.Range("E18").Value = "  +1.17%  "

# Row 19
$ws.Range("D19").Value = "'8.14"
$ws.Range("E19").Value = "  +1.29%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0741"
$ws.Range("E20").Value = "  +1.06%  "

# Row 21
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("D22").Value = "'4.47"
$ws.Range("E22").Value = "  +0.83%  "

# Row 23
$ws.Range("D23").Value = "'9.34"
$ws.Range("E23").Value = "  +1.91%  "

# Row 24
$ws.Range("E24").Value = "  -1.91%  "

# Row 25
$ws.Range("D25").Value = "'146.75"
$ws.Range("E25").Value = "  +0.60%  "

# Row 26
$ws.Range("E26").Value = "  +1.16%  "

# Row 27
$ws.Range("D27").Value = "'16.32"
$ws.Range("E27").Value = "  +2.17%  "

# Row 28
$ws.Range("E28").Value = "  +0.47%  "

# Row 29
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("E30").Value = "  +0.02%  "

# Row 31
$ws.Range("E31").Value = "  +0.25%  "

# Row 32
$ws.Range("B32").Value = "Maker"
$ws.Range("C32").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D32").Value = "1.553.66"
$ws.Range("E32").Value = "  +5.80%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.36"
$ws.Range("E33").Value = "  +0.88%  "

# Row 34
$ws.Range("E34").Value = "  +1.72%  "

# Row 35
$ws.Range("E35").Value = "  +2.33%  "

# Row 36
$ws.Range("E36").Value = "  +3.58%  "

# Row 38
$ws.Range("D38").Value = "'0.925"
$ws.Range("E38").Value = "  +3.17%  "

# Row 39
$ws.Range("E39").Value = "  +1.91%  "

# Row 40
$ws.Range("D40").Value = "'1.06"
$ws.Range("E40").Value = "  -0.05%  "

# Row 41
$ws.Range("E41").Value = "  +0.07%  "

# Row 42
$ws.Range("D42").Value = "'67.79"

# Row 43
$ws.Range("E43").Value = "  -2.98%  "

# Row 44
$ws.Range("D44").Value = "'2.26"
$ws.Range("E44").Value = "  -1.74%  "

# Row 45
$ws.Range("D45").Value = "1.822.97"
$ws.Range("E45").Value = "  +0.69%  "

# Row 46
$ws.Range("D46").Value = "'0.784"
$ws.Range("E46").Value = "  +0.69%  "

# Row 47
$ws.Range("D47").Value = "'90.68"
$ws.Range("E47").Value = "  +0.04%  "

# Row 48
$ws.Range("E48").Value = "  +1.61%  "

# Row 49
$ws.Range("E49").Value = "  +2.59%  "

# Row 50
$ws.Range("E50").Value = "  +2.99%  "

# Row 51
$ws.Range("D51").Value = "'8.04"
$ws.Range("E51").Value = "  +5.03%  "
